$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1540").Value = 'Где на М1 можно найти место для отдыха с рестораном и кафе? @ Отдых'
$ws.Range("C1541").Value = 'Подскажите места для отдыха на М10 возле леса с туалетом и кафе. @ Отдых'
$ws.Range("C1542").Value = 'Где можно остановиться на отдых на М4 перед городом с рестораном и магазином? @ Отдых'
$ws.Range("C1543").Value = 'Есть ли на М5 зоны для отдыха с парком и кафе? @ Отдых'
$ws.Range("C1544").Value = 'Где на М2 можно найти зону отдыха с рестораном и магазином? @ Отдых'
$ws.Range("C1545").Value = 'Какие места для отдыха есть на М6 возле реки с парком? @ Отдых'
$ws.Range("C1546").Value = 'Где можно отдохнуть на М8 перед озером с кафе и рестораном? @ Отдых'
$ws.Range("C1547").Value = 'Есть ли на М11 зоны для отдыха с туалетом и рестораном? @ Отдых'
$ws.Range("C1548").Value = 'Где на М3 можно найти место для отдыха с кафе и рестораном? @ Отдых'
$ws.Range("C1549").Value = 'Подскажите места для отдыха на М7 возле леса с туалетом и магазином. @ Отдых'
$ws.Range("C1550").Value = 'Где можно остановиться на отдых на М1 перед городом с рестораном и кафе? @ Отдых'
$ws.Range("C1551").Value = 'Есть ли на М10 зоны для отдыха с кафе и туалетом? @ Отдых'
$ws.Range("C1552").Value = 'Где на М4 можно найти зону отдыха с рестораном и магазином? @ Отдых'
$ws.Range("C1553").Value = 'Какие места для отдыха есть на М5 возле реки с кафе? @ Отдых'
$ws.Range("C1554").Value = 'Где можно отдохнуть на М2 перед озером с рестораном и магазином? @ Отдых'
$ws.Range("C1555").Value = 'Есть ли на М6 зоны для отдыха с парком и туалетом? @ Отдых'
$ws.Range("C1556").Value = 'Где на М8 можно найти место для отдыха с кафе и рестораном? @ Отдых'
$ws.Range("C1557").Value = 'Какие места для отдыха есть на М11 возле озера с туалетом и рестораном? @ Отдых'
$ws.Range("C1558").Value = 'Где можно отдохнуть на М3 перед городом с парком и кафе? @ Отдых'
$ws.Range("C1559").Value = 'Есть ли на М7 зоны для отдыха с магазином и туалетом? @ Отдых'
$ws.Range("C1560").Value = 'Где на М1 можно найти место для отдыха с рестораном и кафе? @ Отдых'
$ws.Range("C1561").Value = 'Подскажите места для отдыха на М10 возле леса с туалетом и кафе. @ Отдых'
$ws.Range("C1562").Value = 'Где можно остановиться на отдых на М4 перед городом с рестораном и магазином? @ Отдых'
$ws.Range("C1563").Value = 'Есть ли на М5 зоны для отдыха с парком и кафе? @ Отдых'
$ws.Range("C1564").Value = 'Где на М2 можно найти зону отдыха с рестораном и магазином? @ Отдых'
$ws.Range("C1565").Value = 'Какие места для отдыха есть на М6 возле реки с парком? @ Отдых'
$ws.Range("C1566").Value = 'Где можно отдохнуть на М8 перед озером с кафе и рестораном? @ Отдых'
$ws.Range("C1567").Value = 'Есть ли на М11 зоны для отдыха с туалетом и рестораном? @ Отдых'
$ws.Range("C1568").Value = 'Где на М3 можно найти место для отдыха с кафе и рестораном? @ Отдых'
$ws.Range("C1569").Value = 'Подскажите места для отдыха на М7 возле леса с туалетом и магазином. @ Отдых'
$ws.Range("C1570").Value = 'Где на трассе М1 можно найти ресторан с национальной кухней? @ Еда'
$ws.Range("C1572").Value = 'Есть ли кафе на трассе М10, где подают быстрые закуски? @ Еда'
$ws.Range("C1574").Value = 'Какие рестораны на М4 предлагают вегетарианское меню? @ Еда'
$ws.Range("C1576").Value = 'Где на М5 можно найти кафе с домашней выпечкой? @ Еда'
$ws.Range("C1578").Value = 'Можно ли найти на М2 место, где готовят свежие морепродукты? @ Еда'
$ws.Range("C1580").Value = 'Подскажите рестораны на М6, где можно попробовать местные деликатесы? @ Еда'
$ws.Range("C1582").Value = 'Есть ли на М8 кафе с детским меню и игровой зоной? @ Еда'
$ws.Range("C1584").Value = 'Где на М11 можно найти ресторан с видом на озеро? @ Еда'
$ws.Range("C1586").Value = 'Какие рестораны на М3 предлагают бизнес-ланчи? @ Еда'
$ws.Range("C1588").Value = 'Где на М7 можно найти кафе с домашней кухней? @ Еда'
$ws.Range("C1590").Value = 'Есть ли на М1 рестораны с авторской кухней? @ Еда'
$ws.Range("C1592").Value = 'Где на М10 можно найти кафе с быстрым обслуживанием? @ Еда'
$ws.Range("C1594").Value = 'Подскажите рестораны на М4, где можно попробовать местные специализации? @ Еда'
$ws.Range("C1596").Value = 'Где на М5 можно найти кафе с домашними десертами? @ Еда'
$ws.Range("C1598").Value = 'Можно ли найти на М2 место, где готовят свежие овощи и фрукты? @ Еда'
$ws.Range("C1600").Value = 'Есть ли на М6 кафе с традиционными напитками? @ Еда'
$ws.Range("C1602").Value = 'Где на М8 можно найти ресторан с дегустацией местных вин? @ Еда'
$ws.Range("C1604").Value = 'Какие рестораны на М11 предлагают морепродукты? @ Еда'
$ws.Range("C1606").Value = 'Где на М3 можно найти кафе с домашними обедами? @ Еда'
$ws.Range("C1608").Value = 'Подскажите рестораны на М7, где можно попробовать блюда из региональной кухни? @ Еда'
$ws.Range("C1610").Value = 'Есть ли на М1 рестораны с местными деликатесами? @ Еда'
$ws.Range("C1612").Value = 'Где на М10 можно найти кафе с меню для диетического питания? @ Еда'
$ws.Range("C1614").Value = 'Подскажите рестораны на М4, где можно попробовать экзотическую кухню? @ Еда'
$ws.Range("C1616").Value = 'Где на М5 можно найти кафе с домашними пирогами? @ Еда'
$ws.Range("C1618").Value = 'Можно ли найти на М2 место, где готовят свежие морепродукты? @ Еда'
$ws.Range("C1620").Value = 'Есть ли на М6 кафе с традиционными напитками? @ Еда'
$ws.Range("C1622").Value = 'Где на М8 можно найти ресторан с дегустацией местных вин? @ Еда'
$ws.Range("C1624").Value = 'Какие рестораны на М11 предлагают морепродукты? @ Еда'
$ws.Range("C1626").Value = 'Где на М3 можно найти кафе с домашними обедами? @ Еда'
$ws.Range("C1628").Value = 'Подскажите рестораны на М7, где можно попробовать блюда из региональной кухни? @ Еда'

$ws.Range("C1629").Select()
